# Apply the "Office Theme" colour palette to the presentation, replacing
# the current "Integral" theme colours (Design/Theme gallery change).
#
# PowerPoint's MsoThemeColorSchemeIndex order (1-based) used by
# ThemeColorScheme.Item(i).RGB is:
#   1 Dark1, 2 Light1, 3 Dark2, 4 Light2,
#   5-10 Accent1..Accent6, 11 Hyperlink, 12 FollowedHyperlink

$p = $ppt.ActivePresentation

function HexToRgb($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Office Theme colour scheme (the theme being applied).
$officeThemeColors = @(
    "000000",  # Dark1
    "FFFFFF",  # Light1
    "44546A",  # Dark2
    "E7E6E6",  # Light2
    "5B9BD5",  # Accent1
    "ED7D31",  # Accent2
    "A5A5A5",  # Accent3
    "FFC000",  # Accent4
    "4472C4",  # Accent5
    "70AD47",  # Accent6
    "0563C1",  # Hyperlink
    "954F72"   # FollowedHyperlink
)

# The theme is shared by the whole deck (single slide master), so updating
# it through the first slide recolours the master theme used everywhere.
$slide = $p.Slides.Item(1)
$themeColors = $slide.ThemeColorScheme

for ($i = 0; $i -lt $officeThemeColors.Length; $i++) {
    $themeColors.Item($i + 1).RGB = HexToRgb($officeThemeColors[$i])
}
